# ---------------------------------------------------------------------------
# "added nursing demand full functionality"
#
# 1. Nursing Services sheet: add Selected Driver / Technician % /
#    Registered Nurse % / APRN % columns (C:F) with a constant driver
#    ("Patients") and a 30/50/20 split for every specialty row.
# 2. Demand Data by Speciality: rename the "Driver" column header to
#    "Selected Driver", turn on AutoFilter for the full data range and
#    fix up the _FilterDatabase defined name to match.
# 3. Scenarios: tweak the High Scenario percentile value (0.75 -> 0.7).
# 4. Assorted selection / zoom bookkeeping so the saved file re-opens with
#    the same cursor positions as the authored commit.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Nursing Services ("Nursing Services" / sheet1)
# ---------------------------------------------------------------------------
$nursing = $wb.Worksheets.Item("Nursing Services")

$nursing.Range("C1").Value = "Selected Driver"
$nursing.Range("D1").Value = "Technician %"
$nursing.Range("E1").Value = "Registered Nurse %"
$nursing.Range("F1").Value = "APRN %"
$nursing.Range("A1:F1").Font.Bold = $true
$nursing.Range("A1:F1").HorizontalAlignment = -4108

$nursing.Range("C2:C43").HorizontalAlignment = -4108
$nursing.Range("C2:C43").Value = "Patients"

$nursing.Range("D2:D43").HorizontalAlignment = -4108
$nursing.Range("D2:D43").NumberFormat = "0%"
$nursing.Range("D2:D43").Value = 0.3

$nursing.Range("E2:E43").HorizontalAlignment = -4108
$nursing.Range("E2:E43").NumberFormat = "0%"
$nursing.Range("E2:E43").Value = 0.5

$nursing.Range("F2:F43").HorizontalAlignment = -4108
$nursing.Range("F2:F43").NumberFormat = "0%"
$nursing.Range("F2:F43").Value = 0.2

$nursing.Columns.Item(3).ColumnWidth = 18.666666666666668
$nursing.Columns.Item(4).ColumnWidth = 14.0
$nursing.Columns.Item(5).ColumnWidth = 20.0
$nursing.Columns.Item(6).ColumnWidth = 19.333333333333332

# ---------------------------------------------------------------------------
# 2. Health Clusters - cursor / zoom only
# ---------------------------------------------------------------------------
$clusters = $wb.Worksheets.Item("Health Clusters")
$clusters.Activate()
$clusters.Range("A2").Select()
$excel.ActiveWindow.Zoom = 144

# ---------------------------------------------------------------------------
# 3. Benchmarks - cursor only
# ---------------------------------------------------------------------------
$benchmarks = $wb.Worksheets.Item("Benchmarks")
$benchmarks.Activate()
$benchmarks.Range("D28").Select()

# ---------------------------------------------------------------------------
# 4. Demand Data by Speciality
# ---------------------------------------------------------------------------
$demand = $wb.Worksheets.Item("Demand Data by Speciality")
$demand.Activate()

$demand.Range("C1").Value = "Selected Driver"

$demand.Range("A1:W127").AutoFilter()

$filterName = $wb.Names.Item("Demand Data by Speciality!_FilterDatabase")
$filterName.RefersTo = "='Demand Data by Speciality'!`$A`$1:`$W`$127"

$demand.Range("C1").Select()

# ---------------------------------------------------------------------------
# 5. Scenarios - High Scenario percentile 0.75 -> 0.7, cursor
# ---------------------------------------------------------------------------
$scenarios = $wb.Worksheets.Item("Scenarios")
$scenarios.Activate()
$scenarios.Range("C4").Value = 0.7
$scenarios.Range("A6").Select()

# ---------------------------------------------------------------------------
# Finally re-activate Nursing Services so it is the tab that is active when
# the workbook is re-opened (matches the dropped activeTab="3").
# ---------------------------------------------------------------------------
$nursing.Activate()
$nursing.Range("F2").Select()

Write-Output "edit applied"
